# Update ligand/receptor shared-string reorder + refreshed TPM values for rows 2-13,
# and append new "Resolving-Mac" sending-cluster rows 14-17 (new TPM computation).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = @("ECs", "Tgfb2", "Tgfbr1", "ECs", 3, 1, 2.564153333333334, 7.692460000000001, 0.05249149542937438, 0.05249149542937438, 3, 1, 14.861848, 44.585544, 0.09055189482833943, 0.09055189482833945, 38.10805708869334, 342.97251379824, 0.004753204373502969, 0.004753204373502969)
for ($i = 0; $i -lt $row2.Length; $i++) { $ws.Cells.Item(2, $i+1).Value = $row2[$i] }

$row3 = @("ECs", "Tgfb2", "Tgfbr1", "FAPs", 3, 1, 2.564153333333334, 7.692460000000001, 0.05249149542937438, 0.05249149542937438, 3, 1, 31.07813833333333, 93.234415, 0.1893562842131466, 0.1893562842131466, 79.68911200121111, 717.2020080109, 0.0099395945272977, 0.009939594527297702)
for ($i = 0; $i -lt $row3.Length; $i++) { $ws.Cells.Item(3, $i+1).Value = $row3[$i] }

$row4 = @("ECs", "Tgfb2", "Tgfbr1", "MuSCs", 3, 1, 2.564153333333334, 7.692460000000001, 0.05249149542937438, 0.05249149542937438, 3, 1, 18.10188466666667, 54.305654, 0.1102931450066459, 0.1102931450066459, 46.41600790764889, 417.74407116884, 0.005789452117007678, 0.00578945211700768)
for ($i = 0; $i -lt $row4.Length; $i++) { $ws.Cells.Item(4, $i+1).Value = $row4[$i] }

$row5 = @("ECs", "Tgfb2", "Tgfbr1", "Resolving-Mac", 3, 1, 2.564153333333334, 7.692460000000001, 0.05249149542937438, 0.05249149542937438, 3, 1, 100.0833306666667, 300.249992, 0.609798675951868, 0.6097986759518681, 256.6290059400355, 2309.66105346032, 0.03200924441156603, 0.03200924441156604)
for ($i = 0; $i -lt $row5.Length; $i++) { $ws.Cells.Item(5, $i+1).Value = $row5[$i] }

$row6 = @("FAPs", "Tgfb2", "Tgfbr1", "ECs", 3, 1, 16.59481266666667, 49.78443799999999, 0.3397170215679993, 0.3397170215679993, 3, 1, 14.861848, 44.585544, 0.09055189482833943, 0.09055189482833945, 246.6295834404746, 2219.666250964272, 0.03076202000842219, 0.03076202000842219)
for ($i = 0; $i -lt $row6.Length; $i++) { $ws.Cells.Item(6, $i+1).Value = $row6[$i] }

$row7 = @("FAPs", "Tgfb2", "Tgfbr1", "FAPs", 3, 1, 16.59481266666667, 49.78443799999999, 0.3397170215679993, 0.3397170215679993, 3, 1, 31.07813833333333, 93.234415, 0.1893562842131466, 0.1893562842131466, 515.7358836704188, 4641.62295303377, 0.06432755288807372, 0.06432755288807372)
for ($i = 0; $i -lt $row7.Length; $i++) { $ws.Cells.Item(7, $i+1).Value = $row7[$i] }

$row8 = @("FAPs", "Tgfb2", "Tgfbr1", "MuSCs", 3, 1, 16.59481266666667, 49.78443799999999, 0.3397170215679993, 0.3397170215679993, 3, 1, 18.10188466666667, 54.305654, 0.1102931450066459, 0.1102931450066459, 300.3973849569391, 2703.576464612452, 0.0374684587210252, 0.0374684587210252)
for ($i = 0; $i -lt $row8.Length; $i++) { $ws.Cells.Item(8, $i+1).Value = $row8[$i] }

$row9 = @("FAPs", "Tgfb2", "Tgfbr1", "Resolving-Mac", 3, 1, 16.59481266666667, 49.78443799999999, 0.3397170215679993, 0.3397170215679993, 3, 1, 100.0833306666667, 300.249992, 0.609798675951868, 0.6097986759518681, 1660.864123469388, 14947.77711122449, 0.2071589899504782, 0.2071589899504782)
for ($i = 0; $i -lt $row9.Length; $i++) { $ws.Cells.Item(9, $i+1).Value = $row9[$i] }

$row10 = @("MuSCs", "Tgfb2", "Tgfbr1", "ECs", 3, 1, 29.664466, 88.993398, 0.6072695268303631, 0.6072695268303631, 3, 1, 14.861848, 44.585544, 0.09055189482833943, 0.09055189482833945, 440.868784693168, 3967.819062238512, 0.05498940632599849, 0.0549894063259985)
for ($i = 0; $i -lt $row10.Length; $i++) { $ws.Cells.Item(10, $i+1).Value = $row10[$i] }

$row11 = @("MuSCs", "Tgfb2", "Tgfbr1", "FAPs", 3, 1, 29.664466, 88.993398, 0.6072695268303631, 0.6072695268303631, 3, 1, 31.07813833333333, 93.234415, 0.1893562842131466, 0.1893562842131466, 921.9163779324633, 8297.24740139217, 0.1149903011164733, 0.1149903011164733)
for ($i = 0; $i -lt $row11.Length; $i++) { $ws.Cells.Item(11, $i+1).Value = $row11[$i] }

$row12 = @("MuSCs", "Tgfb2", "Tgfbr1", "MuSCs", 3, 1, 29.664466, 88.993398, 0.6072695268303631, 0.6072695268303631, 3, 1, 18.10188466666667, 54.305654, 0.1102931450066459, 0.1102931450066459, 536.9827422302546, 4832.844680072292, 0.06697766598081847, 0.06697766598081849)
for ($i = 0; $i -lt $row12.Length; $i++) { $ws.Cells.Item(12, $i+1).Value = $row12[$i] }

$row13 = @("MuSCs", "Tgfb2", "Tgfbr1", "Resolving-Mac", 3, 1, 29.664466, 88.993398, 0.6072695268303631, 0.6072695268303631, 3, 1, 100.0833306666667, 300.249992, 0.609798675951868, 0.6097986759518681, 2968.91855972809, 26720.26703755281, 0.3703121534070727, 0.3703121534070729)
for ($i = 0; $i -lt $row13.Length; $i++) { $ws.Cells.Item(13, $i+1).Value = $row13[$i] }

$row14 = @("Resolving-Mac", "Tgfb2", "Tgfbr1", "ECs", 1, 0.3333333333333333, 0.025497, 0.076491, 0.000521956172263265, 0.000521956172263265, 3, 1, 14.861848, 44.585544, 0.09055189482833943, 0.09055189482833945, 0.3789325384560001, 3.410392846104, 0.0000472641204157858, 0.0000472641204157858)
for ($i = 0; $i -lt $row14.Length; $i++) { $ws.Cells.Item(14, $i+1).Value = $row14[$i] }

$row15 = @("Resolving-Mac", "Tgfb2", "Tgfbr1", "FAPs", 1, 0.3333333333333333, 0.025497, 0.076491, 0.000521956172263265, 0.000521956172263265, 3, 1, 31.07813833333333, 93.234415, 0.1893562842131466, 0.1893562842131466, 0.792399293085, 7.131593637765, 0.0000988356813018889, 0.00009883568130188891)
for ($i = 0; $i -lt $row15.Length; $i++) { $ws.Cells.Item(15, $i+1).Value = $row15[$i] }

$row16 = @("Resolving-Mac", "Tgfb2", "Tgfbr1", "MuSCs", 1, 0.3333333333333333, 0.025497, 0.076491, 0.000521956172263265, 0.000521956172263265, 3, 1, 18.10188466666667, 54.305654, 0.1102931450066459, 0.1102931450066459, 0.4615437533460001, 4.153893780114, 0.00005756818779454613, 0.00005756818779454614)
for ($i = 0; $i -lt $row16.Length; $i++) { $ws.Cells.Item(16, $i+1).Value = $row16[$i] }

$row17 = @("Resolving-Mac", "Tgfb2", "Tgfbr1", "Resolving-Mac", 1, 0.3333333333333333, 0.025497, 0.076491, 0.000521956172263265, 0.000521956172263265, 3, 1, 100.0833306666667, 300.249992, 0.609798675951868, 0.6097986759518681, 2.551824682008, 22.966422138072, 0.0003182881827510441, 0.0003182881827510442)
for ($i = 0; $i -lt $row17.Length; $i++) { $ws.Cells.Item(17, $i+1).Value = $row17[$i] }

